$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the confidence-interval text strings (row 6 and row 12)
$ws.Range("B6").Value = "(-0.083, 1.082)"
$ws.Range("C6").Value = "(1.457, 3.912)"
$ws.Range("D6").Value = "(11.535, 54.613)"

$ws.Range("B12").Value = "(-0.65, 0.375)"
$ws.Range("C12").Value = "(-0.658, 0.239)"
$ws.Range("D12").Value = "(1.086, 13.933)"

# Update the numeric p-value cells (row 5 and row 11)
$ws.Range("B5").Value = 0.121
$ws.Range("C5").Value = 0.003
$ws.Range("D5").Value = 0

$ws.Range("B11").Value = 0.007
$ws.Range("C11").Value = 0.002

$wb.Save()
